$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal-looking strings (e.g. "9.27") need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers (losing
# the original text formatting, e.g. "0.999"/"580.00" would become 0.999/580 as floats).
$numericLookingCells = @("D5","D6","D12","D14","D17","D20","D22","D23","D24","D26","D28","D30","D33","D34","D37","D38","D39","D43","D45","D46","D48","D49")
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.921.83"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "3.386.10"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "563.74"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "176.46"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "3.381.32"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "53.93"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "9.27"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "3.928.74"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "18.21"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "3.379.26"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "65.900.11"
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").Value = "11.91"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "463.20"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").Value = "4.93"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("D24").Value = "14.85"
$ws.Range("E24").Value = "  +10.03%  "
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").Value = "89.45"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "10.66"
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("D30").Value = "31.11"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "580.00"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").Value = "62.44"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.143"
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "3.59"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").Value = "36.06"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "3.103.41"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "2.85"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.20"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.134"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "140.81"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("E50").Value = "  +9.55%  "
$ws.Range("E51").Value = "  +9.57%  "

# Revert the temporary Text number format back to the default style so the saved
# style indices match the original workbook (only the cell *values* should differ).
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).Style = "Normal"
}
